$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Datetime text updates (stored as plain text strings, not Excel dates) ---
$wsOverview.Range("G2").Value = "2016-10-14 08:25:55"
$wsDeDe.Range("H2").Value = "2016-10-14 08:25:55"
$wsZhCn.Range("H2").Value = "2016-10-14 08:25:45"

# --- Column width changes ---
# Target stored width is 17.2159881591797 characters. The COM layer here
# quantizes ColumnWidth to whole pixels (1/6 character steps with MDW=6),
# so feed it the input value whose rounding lands nearest that target.
$newWidth = 16.333333333333332
$wsOverview.Range("E1").EntireColumn.ColumnWidth = $newWidth
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $newWidth
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $newWidth
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $newWidth
